$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LeaveEntitlement")

# --- Ensure new shared strings are created in the exact order required ---
# Target new shared-string table entries (appended after existing 0-22):
#   23 = "2016-01-01 - 2016-12-31"
#   24 = "Russel"
#   25 = "Thomas"
# Writing H2:H5 in one shot creates a single new shared string (index 23).
# Then writing F3 creates "Russel" (index 24), then F2 creates "Thomas" (index 25).

$ws.Range("H2:H5").Value = "2016-01-01 - 2016-12-31"
$ws.Range("F3").Value = "Russel"
$ws.Range("F2").Value = "Thomas"
$ws.Range("F5").Value = "Russel"

# --- Row 2: fill in remaining (unchanged) values explicitly to be safe ---
$ws.Range("A2").Value = "Y"
$ws.Range("B2").Value = "TC001"
$ws.Range("C2").Value = "N"
$ws.Range("G2").Value = "Maternity US"
$ws.Range("I2").Value = 12

# --- Row 3 (new) ---
$ws.Range("A3").Value = "Y"
$ws.Range("B3").Value = "TC002"
$ws.Range("C3").Value = "N"
$ws.Range("G3").Value = "Vacation US"
$ws.Range("I3").Value = 4

# --- Row 4 (new) ---
$ws.Range("A4").Value = "Y"
$ws.Range("B4").Value = "TC003"
$ws.Range("C4").Value = "Y"
$ws.Range("D4").Value = "United States"
$ws.Range("E4").Value = "Sales"
$ws.Range("G4").Value = "Vacation US"
$ws.Range("I4").Value = 3

# --- Row 5 (new) ---
$ws.Range("A5").Value = "Y"
$ws.Range("B5").Value = "TC004"
$ws.Range("C5").Value = "N"
$ws.Range("G5").Value = "Vacation US"
$ws.Range("I5").Value = 1

# --- Update selection to match target view state ---
$ws.Range("A2:A4").Select() | Out-Null

Write-Host "Edit complete"
